# repull data, push all data, mean calculation
# Update column F ("dSF") values for the rows that changed after repulling data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 4
    6  = -2
    7  = 4
    11 = 0
    12 = -4
    20 = -4
    22 = 3
    23 = -1
    29 = -6
    30 = -6
    36 = -1
    39 = 1
    45 = -2
    46 = 4
    47 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
